# Natmi following Dr Hou advice
# Rewrites the LR-pair result table: recomputed statistics for the
# Wnt11 -> Fzd8 edge across the full ECs / FAPs / sCs cluster grid
# (previously only ECs/FAPs were paired; now all 3x3 combinations are present).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T the 16 numeric statistic columns (same order as the header row).
# NOTE: PowerShell hashtable keys are case-insensitive, so the row-number key is
# named "rr" (not "r") to avoid colliding with the data column "R".
$rows = @(
    @{ rr=2;  A="ECs";  D="ECs";  E=1; F=0.3333333333333333; G=0.062425;            H=0.187275;          I=0.00296197839880675; J=0.00296197839880675; K=2; L=0.6666666666666666; M=1.027778333333333; N=3.083335;  O=0.08020467841353289; P=0.08020467841353289; Q=0.06415906245833333; R=0.577431562125;    S=0.0002375645249441264; T=0.0002375645249441265 }
    @{ rr=3;  A="ECs";  D="FAPs"; E=1; F=0.3333333333333333; G=0.062425;            H=0.187275;          I=0.00296197839880675; J=0.00296197839880675; K=3; L=1;                  M=7.273511666666667; N=21.820535;  O=0.567602609669802;   P=0.567602609669802;   Q=0.4540489657916667;  R=4.086440692125;    S=0.001681226668948293;  T=0.001681226668948293  }
    @{ rr=4;  A="ECs";  D="sCs";  E=1; F=0.3333333333333333; G=0.062425;            H=0.187275;          I=0.00296197839880675; J=0.00296197839880675; K=3; L=1;                  M=4.513153666666667; N=13.539461;  O=0.3521927119166651;  P=0.3521927119166651;  Q=0.2817336176416667;  R=2.535602558775;    S=0.001043187204914331;  T=0.001043187204914331  }
    @{ rr=5;  A="FAPs"; D="ECs";  E=3; F=1;                  G=20.223983;           H=60.671949;         I=0.9595995319797346;  J=0.9595995319797347;  K=2; L=0.6666666666666666; M=1.027778333333333; N=3.083335;  O=0.08020467841353289; P=0.08020467841353289; Q=20.78577154110166;   R=187.071943869915;  S=0.07696437186821128;   T=0.07696437186821128   }
    @{ rr=6;  A="FAPs"; D="FAPs"; E=3; F=1;                  G=20.223983;           H=60.671949;         I=0.9595995319797346;  J=0.9595995319797347;  K=3; L=1;                  M=7.273511666666667; N=21.820535;  O=0.567602609669802;   P=0.567602609669802;   Q=147.0993762969683;   R=1323.894386672715; S=0.544671198589618;     T=0.5446711985896181     }
    @{ rr=7;  A="FAPs"; D="sCs";  E=3; F=1;                  G=20.223983;           H=60.671949;         I=0.9595995319797346;  J=0.9595995319797347;  K=3; L=1;                  M=4.513153666666667; N=13.539461;  O=0.3521927119166651;  P=0.3521927119166651;  Q=91.27394303105433;   R=821.4654872794889; S=0.3379639615219053;    T=0.3379639615219053    }
    @{ rr=8;  A="sCs";  D="ECs";  E=3; F=1;                  G=0.7890326666666666;  H=2.367098;          I=0.0374384896214586;  J=0.03743848962145861; K=2; L=0.6666666666666666; M=1.027778333333333; N=3.083335;  O=0.08020467841353289; P=0.08020467841353289; Q=0.8109506790922221;  R=7.29855611183;     S=0.003002742020377475;  T=0.003002742020377476  }
    @{ rr=9;  A="sCs";  D="FAPs"; E=3; F=1;                  G=0.7890326666666666;  H=2.367098;          I=0.0374384896214586;  J=0.03743848962145861; K=3; L=1;                  M=7.273511666666667; N=21.820535;  O=0.567602609669802;   P=0.567602609669802;   Q=5.739038306381111;   R=51.65134475743;    S=0.0212501844112357;    T=0.0212501844112357    }
    @{ rr=10; A="sCs";  D="sCs";  E=3; F=1;                  G=0.7890326666666666;  H=2.367098;          I=0.0374384896214586;  J=0.03743848962145861; K=3; L=1;                  M=4.513153666666667; N=13.539461;  O=0.3521927119166651;  P=0.3521927119166651;  Q=3.561025672686444;   R=32.04923105417799; S=0.01318556318984542;   T=0.01318556318984542   }
)

$numericCols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $rows) {
    $rowNum = $row.rr
    $ws.Cells.Item($rowNum, 1).Value = $row.A       # Sending cluster
    $ws.Cells.Item($rowNum, 2).Value = "Wnt11"      # Ligand symbol
    $ws.Cells.Item($rowNum, 3).Value = "Fzd8"       # Receptor symbol
    $ws.Cells.Item($rowNum, 4).Value = $row.D       # Target cluster

    for ($i = 0; $i -lt $numericCols.Length; $i++) {
        $col = $numericCols[$i]
        $ws.Cells.Item($rowNum, 5 + $i).Value = $row[$col]
    }
}

Write-Host "Updated rows 2-10 of Sheet1 with recomputed Wnt11/Fzd8 edge statistics"
